$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.819621801376343
$ws.Range("B1").Value = 3.14648699760437
$ws.Range("C1").Value = 1.767034411430359
$ws.Range("D1").Value = 1.416489243507385
$ws.Range("E1").Value = 1.313721060752869
